# Auto-generated edit script applying the Midgardsormr_Profits market-price refresh
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(7, 8).Value = 2000
$ws.Cells.Item(7, 9).Value = 2000
$ws.Cells.Item(7, 11).Value = 2000
$ws.Cells.Item(7, 13).Value = -1888
$ws.Cells.Item(14, 8).Value = 2000
$ws.Cells.Item(14, 9).Value = 2000
$ws.Cells.Item(14, 11).Value = 2000
$ws.Cells.Item(14, 13).Value = -1809
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 14).Value = $null
$ws.Cells.Item(19, 8).Value = 1871.4286
$ws.Cells.Item(19, 9).Value = 1480.2727
$ws.Cells.Item(19, 10).Value = 2301.7
$ws.Cells.Item(19, 11).Value = 1480.2727
$ws.Cells.Item(19, 12).Value = 2301.7
$ws.Cells.Item(19, 13).Value = -1305.2727
$ws.Cells.Item(19, 14).Value = -2651.7
$ws.Cells.Item(40, 8).Value = 3496.7
$ws.Cells.Item(40, 9).Value = 2370.875
$ws.Cells.Item(40, 10).Value = 8000
$ws.Cells.Item(40, 11).Value = 2370.875
$ws.Cells.Item(40, 12).Value = 8000
$ws.Cells.Item(40, 13).Value = -2195.875
$ws.Cells.Item(40, 14).Value = -8350
$ws.Cells.Item(74, 8).Value = 5634.6665
$ws.Cells.Item(74, 9).Value = 5584.75
$ws.Cells.Item(74, 11).Value = 5584.75
$ws.Cells.Item(74, 13).Value = -4648.75
$ws.Cells.Item(77, 8).Value = 5634.6665
$ws.Cells.Item(77, 9).Value = 5584.75
$ws.Cells.Item(77, 11).Value = 27923.75
$ws.Cells.Item(77, 13).Value = -23243.75
$ws.Cells.Item(132, 8).Value = 5104028.5
$ws.Cells.Item(132, 9).Value = 6212844
$ws.Cells.Item(132, 11).Value = 18638532
$ws.Cells.Item(132, 13).Value = -18636002
$ws.Cells.Item(138, 8).Value = 2053673.4
$ws.Cells.Item(138, 9).Value = 2929.6155
$ws.Cells.Item(138, 10).Value = 3120060
$ws.Cells.Item(138, 11).Value = 8788.8465
$ws.Cells.Item(138, 12).Value = 9360180
$ws.Cells.Item(138, 13).Value = -3648.8465
$ws.Cells.Item(138, 14).Value = -9370460

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 16305.972
$ws.Cells.Item(32, 9).Value = 16641.045
$ws.Cells.Item(32, 11).Value = 16641.045
$ws.Cells.Item(32, 13).Value = -16354.045
$ws.Cells.Item(45, 8).Value = 4335.643
$ws.Cells.Item(45, 9).Value = 3278.7144
$ws.Cells.Item(45, 10).Value = 5392.5713
$ws.Cells.Item(45, 11).Value = 3278.7144
$ws.Cells.Item(45, 12).Value = 5392.5713
$ws.Cells.Item(45, 13).Value = -2901.7144
$ws.Cells.Item(45, 14).Value = -6146.5713
$ws.Cells.Item(46, 8).Value = 6087.375
$ws.Cells.Item(46, 9).Value = 5859.4
$ws.Cells.Item(46, 11).Value = 5859.4
$ws.Cells.Item(46, 13).Value = -5540.4
$ws.Cells.Item(61, 8).Value = 9017.16
$ws.Cells.Item(61, 9).Value = 5559.4736
$ws.Cells.Item(61, 11).Value = 5559.4736
$ws.Cells.Item(61, 13).Value = -5347.4736
$ws.Cells.Item(97, 8).Value = 683.7727
$ws.Cells.Item(97, 9).Value = 482.94446
$ws.Cells.Item(97, 11).Value = 482.94446
$ws.Cells.Item(97, 13).Value = 13.05554000000001
$ws.Cells.Item(132, 8).Value = 2001.48
$ws.Cells.Item(132, 9).Value = 1740.9565
$ws.Cells.Item(132, 10).Value = 4997.5
$ws.Cells.Item(132, 11).Value = 5222.8695
$ws.Cells.Item(132, 12).Value = 14992.5
$ws.Cells.Item(132, 13).Value = -2692.8695
$ws.Cells.Item(132, 14).Value = -20052.5
$ws.Cells.Item(136, 8).Value = 9017.16
$ws.Cells.Item(136, 9).Value = 5559.4736
$ws.Cells.Item(136, 11).Value = 16678.4208
$ws.Cells.Item(136, 13).Value = -14128.4208

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2781.4614
$ws.Cells.Item(99, 9).Value = 2803.7273
$ws.Cells.Item(99, 10).Value = 2659
$ws.Cells.Item(99, 11).Value = 2803.7273
$ws.Cells.Item(99, 12).Value = 2659
$ws.Cells.Item(99, 13).Value = -1305.7273
$ws.Cells.Item(99, 14).Value = -5655
$ws.Cells.Item(105, 8).Value = 5210.857
$ws.Cells.Item(105, 10).Value = 4478.8
$ws.Cells.Item(105, 12).Value = 4478.8
$ws.Cells.Item(105, 14).Value = -7972.8
$ws.Cells.Item(134, 8).Value = 3026.4634
$ws.Cells.Item(134, 9).Value = 2813.2163
$ws.Cells.Item(134, 11).Value = 8439.6489
$ws.Cells.Item(134, 13).Value = -5904.6489

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1640.7142
$ws.Cells.Item(58, 9).Value = 1037.7142
$ws.Cells.Item(58, 10).Value = 1942.2142
$ws.Cells.Item(58, 11).Value = 1037.7142
$ws.Cells.Item(58, 12).Value = 1942.2142
$ws.Cells.Item(58, 13).Value = -834.7141999999999
$ws.Cells.Item(58, 14).Value = -2348.2142
$ws.Cells.Item(107, 8).Value = 415.8846
$ws.Cells.Item(107, 9).Value = 251
$ws.Cells.Item(107, 11).Value = 251
$ws.Cells.Item(107, 13).Value = 1669
$ws.Cells.Item(132, 8).Value = 1998.6296
$ws.Cells.Item(132, 9).Value = 1544.091
$ws.Cells.Item(132, 10).Value = 3998.6
$ws.Cells.Item(132, 11).Value = 4632.272999999999
$ws.Cells.Item(132, 12).Value = 11995.8
$ws.Cells.Item(132, 13).Value = -2102.272999999999
$ws.Cells.Item(132, 14).Value = -17055.8
$ws.Cells.Item(134, 8).Value = 3337.3333
$ws.Cells.Item(134, 9).Value = 2506
$ws.Cells.Item(134, 10).Value = 5000
$ws.Cells.Item(134, 11).Value = 7518
$ws.Cells.Item(134, 12).Value = 15000
$ws.Cells.Item(134, 13).Value = -4983
$ws.Cells.Item(134, 14).Value = -20070
$ws.Cells.Item(136, 8).Value = 1640.7142
$ws.Cells.Item(136, 9).Value = 1037.7142
$ws.Cells.Item(136, 10).Value = 1942.2142
$ws.Cells.Item(136, 11).Value = 3113.1426
$ws.Cells.Item(136, 12).Value = 5826.642599999999
$ws.Cells.Item(136, 13).Value = -563.1425999999997
$ws.Cells.Item(136, 14).Value = -10926.6426

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 1338.381
$ws.Cells.Item(132, 9).Value = 1202.2222
$ws.Cells.Item(132, 11).Value = 10819.9998
$ws.Cells.Item(132, 13).Value = -8289.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 7583.2
$ws.Cells.Item(132, 9).Value = 7206.909
$ws.Cells.Item(132, 10).Value = 10342.667
$ws.Cells.Item(132, 11).Value = 21620.727
$ws.Cells.Item(132, 12).Value = 31028.001
$ws.Cells.Item(132, 13).Value = -19090.727
$ws.Cells.Item(132, 14).Value = -36088.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2628.054
$ws.Cells.Item(40, 9).Value = 2521.0857
$ws.Cells.Item(40, 11).Value = 2521.0857
$ws.Cells.Item(40, 13).Value = -2385.0857
$ws.Cells.Item(55, 8).Value = 2882
$ws.Cells.Item(55, 9).Value = 950.5
$ws.Cells.Item(55, 10).Value = 3847.75
$ws.Cells.Item(55, 11).Value = 950.5
$ws.Cells.Item(55, 12).Value = 3847.75
$ws.Cells.Item(55, 13).Value = -777.5
$ws.Cells.Item(55, 14).Value = -4193.75
$ws.Cells.Item(61, 8).Value = 3949.25
$ws.Cells.Item(61, 9).Value = 4999
$ws.Cells.Item(61, 10).Value = 800
$ws.Cells.Item(61, 11).Value = 4999
$ws.Cells.Item(61, 12).Value = 800
$ws.Cells.Item(61, 13).Value = -4797
$ws.Cells.Item(61, 14).Value = -1204
$ws.Cells.Item(100, 8).Value = 3883.4443
$ws.Cells.Item(100, 9).Value = 3101.4
$ws.Cells.Item(100, 10).Value = 4861
$ws.Cells.Item(100, 11).Value = 3101.4
$ws.Cells.Item(100, 12).Value = 4861
$ws.Cells.Item(100, 13).Value = -2560.4
$ws.Cells.Item(100, 14).Value = -5943
$ws.Cells.Item(113, 8).Value = 3949.25
$ws.Cells.Item(113, 9).Value = 4999
$ws.Cells.Item(113, 10).Value = 800
$ws.Cells.Item(113, 11).Value = 4999
$ws.Cells.Item(113, 12).Value = 800
$ws.Cells.Item(113, 13).Value = -2829
$ws.Cells.Item(113, 14).Value = -5140
$ws.Cells.Item(122, 8).Value = 4791
$ws.Cells.Item(122, 9).Value = 4732
$ws.Cells.Item(122, 11).Value = 14196
$ws.Cells.Item(122, 13).Value = -11746
$ws.Cells.Item(132, 8).Value = 4878.5483
$ws.Cells.Item(132, 9).Value = 4813.2856
$ws.Cells.Item(132, 10).Value = 4932.294
$ws.Cells.Item(132, 11).Value = 14439.8568
$ws.Cells.Item(132, 12).Value = 14796.882
$ws.Cells.Item(132, 13).Value = -11909.8568
$ws.Cells.Item(132, 14).Value = -19856.882

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(38, 8).Value = 45452.5
$ws.Cells.Item(38, 10).Value = 45452.5
$ws.Cells.Item(38, 12).Value = 45452.5
$ws.Cells.Item(38, 14).Value = -46398.5
$ws.Cells.Item(49, 8).Value = 20000
$ws.Cells.Item(49, 10).Value = 20000
$ws.Cells.Item(49, 12).Value = 20000
$ws.Cells.Item(49, 14).Value = -20460
$ws.Cells.Item(81, 8).Value = 5025.933
$ws.Cells.Item(81, 9).Value = 4934.4165
$ws.Cells.Item(81, 10).Value = 5392
$ws.Cells.Item(81, 11).Value = 9868.833000000001
$ws.Cells.Item(81, 12).Value = 10784
$ws.Cells.Item(81, 13).Value = -8807.833000000001
$ws.Cells.Item(81, 14).Value = -12906
$ws.Cells.Item(84, 8).Value = 5025.933
$ws.Cells.Item(84, 9).Value = 4934.4165
$ws.Cells.Item(84, 10).Value = 5392
$ws.Cells.Item(84, 11).Value = 49344.165
$ws.Cells.Item(84, 12).Value = 53920
$ws.Cells.Item(84, 13).Value = -44040.165
$ws.Cells.Item(84, 14).Value = -64528
